# Fix the git command shown on the "Get started" slide:
#   "git pull https://github.com/Barstad/Auto-BAHN.git"
# should read
#   "git clone https://github.com/Barstad/Auto-BAHN.git"
#
# The run containing " pull " must become " " followed by a new run
# "clone ", and the run containing "https://" must be split into a
# "https" run and a "://" run (both still carrying the existing
# hyperlink), matching how PowerPoint splits runs when text is edited
# in place.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Locate the paragraph that contains the git command line.
$targetPara = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    if ($para.Text -like "*pull*https*") {
        $targetPara = $para
        break
    }
}

$full = $targetPara.Text

# Replace " pull " with " clone " (keeps it as a single run for now).
$idx = $full.IndexOf(" pull ")
$pullRange = $targetPara.Characters($idx + 1, 6)
$pullRange.Text = " clone "

# Split " clone " into " " and "clone " by touching the first
# character's formatting (no-op re-assignment forces PowerPoint to
# break the run at that boundary without altering the formatting).
$splitA = $targetPara.Characters($idx + 1, 1)
$splitA.Font.Bold = $splitA.Font.Bold

# Re-read the paragraph text/offsets after the edit above.
$full2 = $targetPara.Text
$idx2 = $full2.IndexOf("https://")

# Split "https://" into "https" and "://" the same way.
$splitB = $targetPara.Characters($idx2 + 1, 5)
$splitB.Font.Bold = $splitB.Font.Bold
